# Update NATMI LR-pair metrics (Mfge8-Pdgfrb) with recomputed TPM-based values.
# Ligand (G/H) and Receptor (M/N) expression values changed only for the "ECs"
# cluster; the derived specificity (I/J/O/P) and edge-weight (Q/R/S/T) columns
# are recomputed accordingly for every Sending x Target cluster combination.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 21.64449733333333
$ws.Range("H2").Value = 64.933492
$ws.Range("I2").Value = 0.1098676276771345
$ws.Range("J2").Value = 0.1098676276771345
$ws.Range("M2").Value = 8.488196666666667
$ws.Range("N2").Value = 25.46459
$ws.Range("O2").Value = 0.04138402976425696
$ws.Range("P2").Value = 0.04138402976425696
$ws.Range("Q2").Value = 183.7227501164756
$ws.Range("R2").Value = 1653.50475104828
$ws.Range("S2").Value = 0.004546765173918835
$ws.Range("T2").Value = 0.004546765173918836

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 21.64449733333333
$ws.Range("H3").Value = 64.933492
$ws.Range("I3").Value = 0.1098676276771345
$ws.Range("J3").Value = 0.1098676276771345
$ws.Range("O3").Value = 0.3297460182766552
$ws.Range("P3").Value = 0.3297460182766552
$ws.Range("Q3").Value = 1463.894300841353
$ws.Range("R3").Value = 13175.04870757218
$ws.Range("S3").Value = 0.03622841276403714
$ws.Range("T3").Value = 0.03622841276403714

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 21.64449733333333
$ws.Range("H4").Value = 64.933492
$ws.Range("I4").Value = 0.1098676276771345
$ws.Range("J4").Value = 0.1098676276771345
$ws.Range("O4").Value = 0.6288699519590879
$ws.Range("P4").Value = 0.6288699519590879
$ws.Range("Q4").Value = 2791.843078059268
$ws.Range("R4").Value = 25126.58770253341
$ws.Range("S4").Value = 0.06909244973917852
$ws.Range("T4").Value = 0.06909244973917854

# Row 5: FAPs -> ECs
$ws.Range("H5").Value = 88.285005
$ws.Range("I5").Value = 0.1493784449296822
$ws.Range("J5").Value = 0.1493784449296822
$ws.Range("M5").Value = 8.488196666666667
$ws.Range("N5").Value = 25.46459
$ws.Range("O5").Value = 0.04138402976425696
$ws.Range("P5").Value = 0.04138402976425696
$ws.Range("Q5").Value = 249.79349505255
$ws.Range("R5").Value = 2248.14145547295
$ws.Range("S5").Value = 0.006181882011108386
$ws.Range("T5").Value = 0.006181882011108386

# Row 6: FAPs -> FAPs
$ws.Range("H6").Value = 88.285005
$ws.Range("I6").Value = 0.1493784449296822
$ws.Range("J6").Value = 0.1493784449296822
$ws.Range("O6").Value = 0.3297460182766552
$ws.Range("P6").Value = 0.3297460182766552
$ws.Range("R6").Value = 17913.08622402832
$ws.Range("S6").Value = 0.0492569474319213
$ws.Range("T6").Value = 0.04925694743192131

# Row 7: FAPs -> MuSCs
$ws.Range("H7").Value = 88.285005
$ws.Range("I7").Value = 0.1493784449296822
$ws.Range("J7").Value = 0.1493784449296822
$ws.Range("O7").Value = 0.6288699519590879
$ws.Range("P7").Value = 0.6288699519590879
$ws.Range("Q7").Value = 3795.85130129268
$ws.Range("S7").Value = 0.09393961548665249
$ws.Range("T7").Value = 0.0939396154866525

# Row 8: MuSCs -> ECs
$ws.Range("I8").Value = 0.7407539273931834
$ws.Range("J8").Value = 0.7407539273931834
$ws.Range("M8").Value = 8.488196666666667
$ws.Range("N8").Value = 25.46459
$ws.Range("O8").Value = 0.04138402976425696
$ws.Range("P8").Value = 0.04138402976425696
$ws.Range("Q8").Value = 1238.702897091672
$ws.Range("R8").Value = 11148.32607382505
$ws.Range("S8").Value = 0.03065538257922974
$ws.Range("T8").Value = 0.03065538257922974

# Row 9: MuSCs -> FAPs
$ws.Range("I9").Value = 0.7407539273931834
$ws.Range("J9").Value = 0.7407539273931834
$ws.Range("O9").Value = 0.3297460182766552
$ws.Range("P9").Value = 0.3297460182766552
$ws.Range("S9").Value = 0.2442606580806967
$ws.Range("T9").Value = 0.2442606580806967

# Row 10: MuSCs -> MuSCs
$ws.Range("I10").Value = 0.7407539273931834
$ws.Range("J10").Value = 0.7407539273931834
$ws.Range("O10").Value = 0.6288699519590879
$ws.Range("P10").Value = 0.6288699519590879
$ws.Range("S10").Value = 0.4658378867332569
$ws.Range("T10").Value = 0.4658378867332569
